$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week day block rows 20-25 corresponds to 2018/3/8 (Thursday).
# Today's task is all complete: mark every actual-completion cell "yes"
# and set the overall completion percentage (merged H20:H25) to 100%.
$ws.Range("F20").Value = "yes"
$ws.Range("F21").Value = "yes"
$ws.Range("F22").Value = "yes"
$ws.Range("F23").Value = "yes"
$ws.Range("F24").Value = "yes"
$ws.Range("F25").Value = "yes"

$ws.Range("H20").Value = 1

$ws.Range("L27").Select()
